# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to the newly scraped counts (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 7425   # F2: 7413 -> 7425
$ws1.Cells.Item(2, 7).Value = 60     # G2: 54   -> 60
$ws1.Cells.Item(3, 6).Value = 7402   # F3: 7374 -> 7402
$ws1.Cells.Item(4, 6).Value = 97     # F4: 95   -> 97
$ws1.Cells.Item(5, 6).Value = 184    # F5: 183  -> 184
$ws1.Cells.Item(9, 6).Value = 95     # F9: 94   -> 95
$ws1.Cells.Item(10, 6).Value = 131   # F10: 128 -> 131
$ws1.Cells.Item(11, 6).Value = 213   # F11: 211 -> 213
$ws1.Cells.Item(12, 6).Value = 100   # F12: 97  -> 100
$ws1.Cells.Item(13, 6).Value = 670   # F13: 668 -> 670
$ws1.Cells.Item(14, 6).Value = 537   # F14: 522 -> 537
$ws1.Cells.Item(16, 6).Value = 32    # F16: 30  -> 32
$ws1.Cells.Item(19, 6).Value = 76    # F19: 75  -> 76

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 7425   # F2: 7413 -> 7425
$ws4.Cells.Item(2, 7).Value = 60     # G2: 54   -> 60
$ws4.Cells.Item(3, 6).Value = 7402   # F3: 7374 -> 7402
$ws4.Cells.Item(4, 6).Value = 97     # F4: 95   -> 97
$ws4.Cells.Item(5, 6).Value = 184    # F5: 183  -> 184
$ws4.Cells.Item(9, 6).Value = 95     # F9: 94   -> 95
$ws4.Cells.Item(10, 6).Value = 131   # F10: 128 -> 131
$ws4.Cells.Item(11, 6).Value = 213   # F11: 211 -> 213
$ws4.Cells.Item(12, 6).Value = 100   # F12: 97  -> 100
$ws4.Cells.Item(13, 6).Value = 670   # F13: 668 -> 670
$ws4.Cells.Item(14, 6).Value = 537   # F14: 522 -> 537
$ws4.Cells.Item(16, 6).Value = 32    # F16: 30  -> 32
# NOTE: F19 on "全部类型" is already 76 in the source workbook, so no change here.
